# Delete the last slide (slide12.xml, sldId 267) from the presentation.
$p = $ppt.ActivePresentation
$p.Slides.Item($p.Slides.Count).Delete()
